$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    27.44583381195966,
    30.5458117586968,
    31.60593708186514,
    32.42329172343318,
    27.51499558622146,
    32.71195847311773,
    29.91932121404855,
    28.5680508162996,
    28.54489462734267,
    28.5812396920409,
    22.62771968201722,
    26.51569043355869,
    24.99743035012028,
    18.5536949996409,
    18.51927836362415,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277,
    14.13501495067277
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
